$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(23356, "Natália Farias", "Marketing", "Viagem de negocios", 7, 45080, 9315.950000000001),
    @(81588, "Pedro Miguel Monteiro", "Atendimento ao Cliente", "Outros", 6, 45084, 3259.68),
    @(79661, "Dra. Eloah Lima", "Juridico", "Outros", 8, 45094, 9027.42),
    @(59138, "Maria Helena Brito", "Vendas", "Doenca", 4, 45095, 2605.51),
    @(26737, "Sofia Freitas", "Atendimento ao Cliente", "Doenca", 8, 45098, 2750.93),
    @(49672, "Brenda Sampaio", "Atendimento ao Cliente", "Doenca", 6, 45104, 4339.02),
    @(22060, "Dra. Maria Laura Moraes", "Recursos Humanos", "Doenca", 6, 45087, 3583.23),
    @(83399, "Yago Moraes", "Operacoes", "Consulta medica", 1, 45090, 7444.47),
    @(67308, "Thales Garcia", "Vendas", "Doenca", 8, 45090, 2298.41),
    @(44564, "Maria Flor Fonseca", "P&D", "Problemas pessoais", 4, 45095, 3117.5)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rec[6]
    $rowIndex++
}

$wb.Save()
